# This script updates the Java stack-trace text embedded in the document to
# reflect line-number / method-name changes caused by the POI 3.16 -> 3.17
# upgrade (commit "Fixed #253 Moving from POI 3.16 to 3.17").
#
# The whole stack trace lives inside a single <w:t> run, as literal text
# containing tab ([char]9) and line-feed ([char]10) characters (no <w:tab/>
# or <w:br/> child elements are used for the interior lines). We therefore
# perform three scoped Find & Replace operations, each keyed on a distinct,
# uniquely-identifying multi-line excerpt, rather than rewriting the whole
# run at once.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) caseLet/doSwitch/caseBlock line numbers; caseTemplate renamed to
#    caseDocumentTemplate (with new line numbers); generate/M2DocUtils/
#    AbstractTemplatesTestSuite line numbers; GeneratedMethodAccessor
#    index change. This also removes one duplicated
#    "caseDocumentTemplate(...) / doSwitch(...)" block that is no longer
#    present in the updated stack trace.
# ---------------------------------------------------------------------
$old1 = @(
    ([char]9 + 'at org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)')
    ([char]9 + 'at org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)')
    ([char]9 + 'at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseLet(M2DocEvaluator.java:847)')
    ([char]9 + 'at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseLet(M2DocEvaluator.java:1)')
    ([char]9 + 'at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:275)')
    ([char]9 + 'at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)')
    ([char]9 + 'at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)')
    ([char]9 + 'at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:836)')
    ([char]9 + 'at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1034)')
    ([char]9 + 'at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)')
    ([char]9 + 'at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:183)')
    ([char]9 + 'at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)')
    ([char]9 + 'at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)')
    ([char]9 + 'at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:836)')
    ([char]9 + 'at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseTemplate(M2DocEvaluator.java:297)')
    ([char]9 + 'at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseTemplate(M2DocEvaluator.java:1)')
    ([char]9 + 'at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:201)')
    ([char]9 + 'at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)')
    ([char]9 + 'at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)')
    ([char]9 + 'at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:836)')
    ([char]9 + 'at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:259)')
    ([char]9 + 'at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)')
    ([char]9 + 'at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:246)')
    ([char]9 + 'at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)')
    ([char]9 + 'at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)')
    ([char]9 + 'at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:836)')
    ([char]9 + 'at org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:252)')
    ([char]9 + 'at org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:691)')
    ([char]9 + 'at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:396)')
    ([char]9 + 'at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:318)')
    ([char]9 + 'at sun.reflect.GeneratedMethodAccessor4.invoke(Unknown Source)')
    ([char]9 + 'at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)')
    ([char]9 + 'at java.lang.reflect.Method.invoke(Method.java:498)')
    ([char]9 + 'at org.junit.runners.model.FrameworkMethod$1.runReflectiveCall(FrameworkMethod.java:50)')
    ([char]9 + 'at org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)')
    ([char]9 + 'at org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:47)')
    ([char]9 + 'at org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)')
    ([char]9 + 'at org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:325)')
    ([char]9 + 'at org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:78)')
    ([char]9 + 'at org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:57)')
) -join [char]10

$new1 = @(
    ([char]9 + 'at org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)')
    ([char]9 + 'at org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)')
    ([char]9 + 'at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseLet(M2DocEvaluator.java:1050)')
    ([char]9 + 'at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseLet(M2DocEvaluator.java:1)')
    ([char]9 + 'at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:314)')
    ([char]9 + 'at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)')
    ([char]9 + 'at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)')
    ([char]9 + 'at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)')
    ([char]9 + 'at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1254)')
    ([char]9 + 'at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)')
    ([char]9 + 'at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)')
    ([char]9 + 'at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)')
    ([char]9 + 'at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)')
    ([char]9 + 'at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)')
    ([char]9 + 'at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:275)')
    ([char]9 + 'at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)')
    ([char]9 + 'at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)')
    ([char]9 + 'at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)')
    ([char]9 + 'at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)')
    ([char]9 + 'at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)')
    ([char]9 + 'at org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:264)')
    ([char]9 + 'at org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:712)')
    ([char]9 + 'at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:459)')
    ([char]9 + 'at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:369)')
    ([char]9 + 'at sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)')
    ([char]9 + 'at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)')
    ([char]9 + 'at java.lang.reflect.Method.invoke(Method.java:498)')
    ([char]9 + 'at org.junit.runners.model.FrameworkMethod$1.runReflectiveCall(FrameworkMethod.java:50)')
    ([char]9 + 'at org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)')
    ([char]9 + 'at org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:47)')
    ([char]9 + 'at org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)')
    ([char]9 + 'at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)')
    ([char]9 + 'at org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:325)')
    ([char]9 + 'at org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:78)')
    ([char]9 + 'at org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:57)')
) -join [char]10

$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
if (-not $found1) {
    Write-Output "WARNING: block 1 not found"
}

# ---------------------------------------------------------------------
# 2) Remove the now-absent "RunBefores.evaluate(RunBefores.java:26)"
#    frame that used to precede the second "RunAfters" frame.
# ---------------------------------------------------------------------
$old2 = @(
    ([char]9 + 'at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)')
    ([char]9 + 'at org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)')
    ([char]9 + 'at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)')
    ([char]9 + 'at org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)')
    ([char]9 + 'at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)')
    ([char]9 + 'at org.junit.runners.ParentRunner.run(ParentRunner.java:363)')
    ([char]9 + 'at org.junit.runners.Suite.runChild(Suite.java:128)')
) -join [char]10

$new2 = @(
    ([char]9 + 'at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)')
    ([char]9 + 'at org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)')
    ([char]9 + 'at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)')
    ([char]9 + 'at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)')
    ([char]9 + 'at org.junit.runners.ParentRunner.run(ParentRunner.java:363)')
    ([char]9 + 'at org.junit.runners.Suite.runChild(Suite.java:128)')
) -join [char]10

$found2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
if (-not $found2) {
    Write-Output "WARNING: block 2 not found"
}

# ---------------------------------------------------------------------
# 3) Update the Eclipse JDT RemoteTestRunner line numbers at the bottom
#    of the stack trace.
# ---------------------------------------------------------------------
$old3 = @(
    ([char]9 + 'at org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)')
    ([char]9 + 'at org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)')
    ([char]9 + 'at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:459)')
    ([char]9 + 'at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:675)')
    ([char]9 + 'at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:382)')
    ([char]9 + 'at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:192)')
) -join [char]10

$new3 = @(
    ([char]9 + 'at org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)')
    ([char]9 + 'at org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)')
    ([char]9 + 'at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:539)')
    ([char]9 + 'at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:761)')
    ([char]9 + 'at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:461)')
    ([char]9 + 'at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:207)')
) -join [char]10

$found3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
if (-not $found3) {
    Write-Output "WARNING: block 3 not found"
}

Write-Output "Block1 replaced: $found1"
Write-Output "Block2 replaced: $found2"
Write-Output "Block3 replaced: $found3"
